$d = $word.ActiveDocument

# 1) Add a 5-pt-space paragraph border (top/left/bottom/right) to the first
#    paragraph, and change its left indent from 120 twips (6pt) to 225 twips
#    (11.25pt).
$p1 = $d.Paragraphs(1)
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5
$p1.Format.LeftIndent = 11.25

# 2) Remove the trailing, standalone " " run at the end of the first
#    paragraph (leaving the "**ID__AFFARS_MP5332_4__ID**" run and the
#    paragraph mark intact). The paragraph range's last character is the
#    paragraph mark itself, so the trailing space sits one character before
#    that.
$pEnd = $p1.Range.End
$trailingSpace = $d.Range($pEnd - 2, $pEnd - 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}
